# Update with Correct Forecast output
# - Rename Sheet1 -> "Sales vs PO"
# - Insert an "Order Week" column holding the original weekly dates, shift
#   ds to the week-ending date (+6 days) and reset PO_Requested_Qty to the
#   per-day series (0 except real PO intake moved to the new Weekly Growth
#   sheet).
# - Add "Weekly Growth", "Volume Insights" and "Prediction Info" sheets with
#   the derived forecast analytics.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# --- Sales vs PO: insert the "Order Week" column before PO_Requested_Qty ---
$ws1.Columns.Item(3).Insert()
$ws1.Range("C1").Value = "Order Week"
$ws1.Range("C2:C14").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Original weekly-start dates (previously in column A) now live in column C;
# column A becomes the week-ending ("ds") date, 6 days later.
$origDates = @(45565, 45572, 45579, 45586, 45593, 45600, 45607, 45614, 45621, 45628, 45635, 45642, 45649)

for ($i = 0; $i -lt $origDates.Length; $i++) {
    $row = $i + 2
    $orig = $origDates[$i]
    $ws1.Cells.Item($row, 1).Value = $orig + 6
    $ws1.Cells.Item($row, 3).Value = $orig
    $ws1.Cells.Item($row, 4).Value = 0
}

# --- Weekly Growth sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"
# Re-use the existing header style (bold, centered, thin border) from Sales vs PO.
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

$growthDates = @(45572, 45586, 45593, 45621, 45628, 45635)
$growthQty = @(200, 30, 80, 10, 10, 20)
$growthPct = @(0, -85, 166.6666666666667, -87.5, 0, 100)

for ($i = 0; $i -lt $growthDates.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $growthDates[$i]
    $ws2.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws2.Cells.Item($row, 2).Value = $growthQty[$i]
    $ws2.Cells.Item($row, 3).Value = $growthPct[$i]
}

# --- Volume Insights sheet ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"
$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A2").Value = 350
$ws3.Range("B2").Value = 58.33333333333334
$ws3.Range("C2").Value = 200
$ws3.Range("D2").Value = 10

# --- Prediction Info sheet ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

$ws4.Range("A2").Value = 0
